$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.833.99"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "1.856.81"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5033"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07153"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8908"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "1.846.44"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.219"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008498"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "26.878.48"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.010"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").Value = "2.087.13"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.402"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.793"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.633"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.649"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09221"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05078"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7442"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.971"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.247"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.503"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.092"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.484"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.420"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1461"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4629"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.556"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.83%  "
